# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Fri Oct 25 21:42:55 UTC 2024 with GitHub Actions".
# Only column D (Price) and column E (Volume(1h)) text values change, on
# rows 2-51 of Sheet1; everything else (coin name, link, row count, styles)
# stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, [string]$text) {
    # Writes $text into $range as literal text, even when it looks like a
    # number (e.g. "0.999", "4.11", "26.00") so Excel's input parser does
    # not silently convert it into a numeric cell -- the source cells are
    # plain text/inline strings. The cell's existing style is restored
    # afterwards so no formatting change is left behind.
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") '67.186.85'
$ws.Range("E2").Value = '  -1.86%  '

Set-TextValue $ws.Range("D3") '2.487.08'
$ws.Range("E3").Value = '  -2.22%  '

Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.12%  '

Set-TextValue $ws.Range("D5") '586.20'
$ws.Range("E5").Value = '  -1.63%  '

Set-TextValue $ws.Range("D6") '168.24'
$ws.Range("E6").Value = '  -5.37%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -2.88%  '

Set-TextValue $ws.Range("D9") '2.487.49'
$ws.Range("E9").Value = '  -2.21%  '

$ws.Range("E10").Value = '  -4.31%  '

$ws.Range("E11").Value = '  -0.07%  '

Set-TextValue $ws.Range("D12") '0.340'
$ws.Range("E12").Value = '  -2.85%  '

$ws.Range("E13").Value = '  -4.26%  '

Set-TextValue $ws.Range("D14") '26.00'
$ws.Range("E14").Value = '  -4.05%  '

Set-TextValue $ws.Range("D15") '2.930.36'
$ws.Range("E15").Value = '  -2.67%  '

$ws.Range("E16").Value = '  -3.68%  '

Set-TextValue $ws.Range("D17") '66.768.65'
$ws.Range("E17").Value = '  -2.46%  '

Set-TextValue $ws.Range("D18") '2.477.61'
$ws.Range("E18").Value = '  -3.39%  '

$ws.Range("E19").Value = '  +1.37%  '

$ws.Range("E20").Value = '  -2.95%  '

Set-TextValue $ws.Range("D21") '360.56'
$ws.Range("E21").Value = '  -2.19%  '

Set-TextValue $ws.Range("D22") '4.11'
$ws.Range("E22").Value = '  -2.87%  '

Set-TextValue $ws.Range("D23") '4.43'
$ws.Range("E23").Value = '  -6.46%  '

$ws.Range("E24").Value = '  +0.14%  '

Set-TextValue $ws.Range("D25") '70.83'
$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("E26").Value = '  -6.10%  '

$ws.Range("E27").Value = '  -8.59%  '

Set-TextValue $ws.Range("D28") '0.998'
$ws.Range("E28").Value = '  +0.17%  '

Set-TextValue $ws.Range("D29") '2.612.58'
$ws.Range("E29").Value = '  -2.61%  '

Set-TextValue $ws.Range("D30") '0.0₃0936'
$ws.Range("E30").Value = '  -6.54%  '

Set-TextValue $ws.Range("D31") '8.09'
$ws.Range("E31").Value = '  -2.61%  '

Set-TextValue $ws.Range("D32") '508.66'
$ws.Range("E32").Value = '  -6.62%  '

Set-TextValue $ws.Range("D33") '1.84'
$ws.Range("E33").Value = '  -2.79%  '

$ws.Range("E34").Value = '  -5.64%  '

Set-TextValue $ws.Range("D35") '0.999'
$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("E36").Value = '  -2.46%  '

Set-TextValue $ws.Range("D37") '158.51'
$ws.Range("E37").Value = '  +0.59%  '

Set-TextValue $ws.Range("D38") '19.08'
$ws.Range("E38").Value = '  +0.69%  '

$ws.Range("E39").Value = '  -3.73%  '

$ws.Range("E40").Value = '  -0.77%  '

Set-TextValue $ws.Range("D41") '1.74'
$ws.Range("E41").Value = '  -4.82%  '

Set-TextValue $ws.Range("D42") '4.97'
$ws.Range("E42").Value = '  -5.49%  '

Set-TextValue $ws.Range("D43") '0.336'
$ws.Range("E43").Value = '  -6.10%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("E45").Value = '  -4.02%  '

Set-TextValue $ws.Range("D46") '39.40'
$ws.Range("E46").Value = '  -1.88%  '

$ws.Range("E47").Value = '  -4.21%  '

Set-TextValue $ws.Range("D48") '0.540'
$ws.Range("E48").Value = '  -4.66%  '

Set-TextValue $ws.Range("D49") '3.61'
$ws.Range("E49").Value = '  -3.81%  '

Set-TextValue $ws.Range("D50") '0.0₆0267'
$ws.Range("E50").Value = '  -5.32%  '

$ws.Range("E51").Value = '  -4.18%  '
